$d = $word.ActiveDocument

# 1) Replace the subtitle text (Range.Text avoids Find/Replace smart-quote autocorrect,
#    preserving the literal straight quotes/apostrophes used in the new subtitle)
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = 'Dopo l''intervista a Repubblica del capogruppo di Italia Viva che ipotizza un''unione dei moderati i consiglieri Antonella Leto, Luigi Carollo, Alberto Mangano e Ramon La Torre rispondono: "Aiuti ai cittadini non alle multinazionali dei termovalorizzatori"'

# 2) Delete the page-break paragraph and the five word-definition paragraphs (paragraphs 4..9)
$count = $d.Paragraphs.Count
$pStart = $d.Paragraphs.Item(4)
$pEnd = $d.Paragraphs.Item($count)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

# 3) Rewrite the body paragraph text: one run, two w:t separated by a pair of line breaks
$vtab = [char]11
$bodyText = '“Il ‘progresso’ che annuncia Davide Faraone per Palermo, chiamando a raccolta Forza Italia, Udc e autonomisti e chi altro si riconosce nella prospettiva di privatizzare tutto, compresa l''aria che respiriamo, è diametralmente opposta alla visione di Sinistra Comune”. Dopo l’intervista rilasciata a Repubblica dal capogruppo di Italia viva al Senato, che dopo gli scontri al Comune di Palermo entra a gamba tesa sugli alleati e descrive la sua visione di futuro, nella quale aggregare “i moderati” su temi come la realizzazione di un termovalorizzatore a Palermo e il Ponte sullo Stretto, Sinistra Comune va all’attacco: “La nostra – dicono in una nota Antonella Leto, Luigi Carollo, Alberto Mangano e Ramon La Torre – è una visione di città che guarda all''economia circolare ed alla sostenibilità ambientale, al benessere della comunità, alla programmazione dell''impiantistica necessaria a costruire un futuro green e nuovi posti di lavoro. I rifiuti, se si chiude il ciclo virtuoso, cominciando magari dall''approvare il regolamento che da mesi attende di essere discusso in Consiglio comunale, sono risorse economiche, è denaro contante che può tornare al cittadino anziché far fare profitto a qualche multinazionale della termovalorizzazione”.' + $vtab + $vtab + 'Così, per la sinistra-sinistra finita all’angolo per le polemiche sull’assessore alla Mobilità Giusto Catania, è l’occasione per andare all’attacco dei renziani su Rap, l’azienda che gestisce i rifiuti nel capoluogo e che oggi è sotto i riflettori per l’inchiesta della Procura di Palermo che coinvolge, fra gli altri, un funzionario della società, Vincenzo Bonanno: “Per fare questo – prosegue Sinistra Comune - non basta che Rap sia una spa a totale capitale pubblico, è necessario che sia sotto il controllo e la proposta popolare, che divenga come tutte le altre società partecipate del Comune, un’azienda pubblica e partecipativa. Bisogna che l''interesse della comunità Palermo sia anteposto agli interessi dei privati, non sempre leciti, bisogna che la gestione sia trasparente, efficace ed economica, in grado di lavorare in sinergia con le altre aziende pubbliche per generare risparmi e programmare il futuro sostenibile”.try { MNZ_RICH(''Bottom''); } catch(e) {}'
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = $bodyText
